$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 129.84616
$ws.Range("I6").Value = 129.84616
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 389.53848
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -277.53848
$ws.Range("N6").Value = $null
$ws.Range("H17").Value = 331.66666
$ws.Range("J17").Value = 331.66666
$ws.Range("L17").Value = 994.9999799999999
$ws.Range("N17").Value = -1330.99998
$ws.Range("H98").Value = 1870.3889
$ws.Range("I98").Value = 790
$ws.Range("J98").Value = 7272.3335
$ws.Range("K98").Value = 790
$ws.Range("L98").Value = 7272.3335
$ws.Range("M98").Value = 708
$ws.Range("N98").Value = -10268.3335
$ws.Range("H122").Value = 1870.3889
$ws.Range("I122").Value = 790
$ws.Range("J122").Value = 7272.3335
$ws.Range("K122").Value = 2370
$ws.Range("L122").Value = 21817.0005
$ws.Range("M122").Value = 80
$ws.Range("N122").Value = -26717.0005
$ws.Range("H131").Value = 1494.2174
$ws.Range("I131").Value = 1257.8
$ws.Range("J131").Value = 1937.5
$ws.Range("K131").Value = 3773.4
$ws.Range("L131").Value = 5812.5
$ws.Range("M131").Value = 1266.6
$ws.Range("N131").Value = -15892.5
$ws.Range("H138").Value = 2278.5134
$ws.Range("I138").Value = 2754.182
$ws.Range("J138").Value = 2077.2693
$ws.Range("K138").Value = 8262.545999999998
$ws.Range("L138").Value = 6231.8079
$ws.Range("M138").Value = -3122.545999999998
$ws.Range("N138").Value = -16511.8079

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1073177.6
$ws.Range("I32").Value = 1114451.6
$ws.Range("J32").Value = 41330
$ws.Range("K32").Value = 1114451.6
$ws.Range("L32").Value = 41330
$ws.Range("M32").Value = -1114164.6
$ws.Range("N32").Value = -41904
$ws.Range("H110").Value = 1734
$ws.Range("I110").Value = 1841.55
$ws.Range("J110").Value = 1518.9
$ws.Range("K110").Value = 1841.55
$ws.Range("L110").Value = 1518.9
$ws.Range("M110").Value = 203.45
$ws.Range("N110").Value = -5608.9
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 24890.75
$ws.Range("I82").Value = 10000
$ws.Range("J82").Value = 29854.334
$ws.Range("K82").Value = 10000
$ws.Range("L82").Value = 29854.334
$ws.Range("M82").Value = -9617
$ws.Range("N82").Value = -30620.334
$ws.Range("H85").Value = 24890.75
$ws.Range("I85").Value = 10000
$ws.Range("J85").Value = 29854.334
$ws.Range("K85").Value = 10000
$ws.Range("L85").Value = 29854.334
$ws.Range("M85").Value = -8674
$ws.Range("N85").Value = -32506.334
$ws.Range("H94").Value = 1044.0769
$ws.Range("I94").Value = 914.8182
$ws.Range("J94").Value = 1755
$ws.Range("K94").Value = 914.8182
$ws.Range("L94").Value = 1755
$ws.Range("M94").Value = -463.8182
$ws.Range("N94").Value = -2657

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 3575
$ws.Range("J13").Value = 3575
$ws.Range("L13").Value = 3575
$ws.Range("N13").Value = -3853
$ws.Range("H31").Value = 1183336
$ws.Range("I31").Value = 992.4838999999999
$ws.Range("J31").Value = 1980132.6
$ws.Range("K31").Value = 992.4838999999999
$ws.Range("L31").Value = 1980132.6
$ws.Range("M31").Value = -697.4838999999999
$ws.Range("N31").Value = -1980722.6
$ws.Range("H34").Value = 1183336
$ws.Range("I34").Value = 992.4838999999999
$ws.Range("J34").Value = 1980132.6
$ws.Range("K34").Value = 992.4838999999999
$ws.Range("L34").Value = 1980132.6
$ws.Range("M34").Value = -790.4838999999999
$ws.Range("N34").Value = -1980536.6
$ws.Range("H80").Value = 26666.666
$ws.Range("J80").Value = 26666.666
$ws.Range("L80").Value = 26666.666
$ws.Range("N80").Value = -28912.666
$ws.Range("H83").Value = 26666.666
$ws.Range("J83").Value = 26666.666
$ws.Range("L83").Value = 79999.99800000001
$ws.Range("N83").Value = -91231.99800000001
$ws.Range("H132").Value = 17243608
$ws.Range("I132").Value = 33334706
$ws.Range("J132").Value = 3144.1428
$ws.Range("K132").Value = 100004118
$ws.Range("L132").Value = 9432.428400000001
$ws.Range("M132").Value = -100001588
$ws.Range("N132").Value = -14492.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 928.75
$ws.Range("I11").Value = 928.75
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 2786.25
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -2646.25
$ws.Range("N11").Value = $null
$ws.Range("H20").Value = 3599.75
$ws.Range("J20").Value = 3599.75
$ws.Range("L20").Value = 10799.25
$ws.Range("N20").Value = -11253.25
$ws.Range("H107").Value = 915.8095
$ws.Range("I107").Value = 335
$ws.Range("J107").Value = 1859.625
$ws.Range("K107").Value = 1005
$ws.Range("L107").Value = 5578.875
$ws.Range("M107").Value = 915
$ws.Range("N107").Value = -9418.875
$ws.Range("H118").Value = 2740.5386
$ws.Range("I118").Value = 588.1667
$ws.Range("J118").Value = 4585.4287
$ws.Range("K118").Value = 1764.5001
$ws.Range("L118").Value = 13756.2861
$ws.Range("M118").Value = -521.5001
$ws.Range("N118").Value = -16242.2861
$ws.Range("H131").Value = 1348.95
$ws.Range("I131").Value = 1183.909
$ws.Range("J131").Value = 1411.5518
$ws.Range("K131").Value = 3551.727
$ws.Range("L131").Value = 4234.6554
$ws.Range("M131").Value = 1488.273
$ws.Range("N131").Value = -14314.6554
$ws.Range("H141").Value = 3893
$ws.Range("I141").Value = 3456.7856
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 10370.3568
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -5190.356800000001
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 10
$ws.Range("I4").Value = 10
$ws.Range("K4").Value = 10
$ws.Range("M4").Value = 102
$ws.Range("H5").Value = 10004
$ws.Range("I5").Value = 10004
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 10004
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -9892
$ws.Range("N5").Value = $null
$ws.Range("H113").Value = 33341408
$ws.Range("I113").Value = 83350830
$ws.Range("J113").Value = 1789.1111
$ws.Range("K113").Value = 83350830
$ws.Range("L113").Value = 1789.1111
$ws.Range("M113").Value = -83348660
$ws.Range("N113").Value = -6129.1111
$ws.Range("H122").Value = 1051.1666
$ws.Range("I122").Value = 926.5
$ws.Range("J122").Value = 1086.7858
$ws.Range("K122").Value = 2779.5
$ws.Range("L122").Value = 3260.3574
$ws.Range("M122").Value = -329.5
$ws.Range("N122").Value = -8160.357400000001
$ws.Range("H141").Value = 22500
$ws.Range("J141").Value = 22500
$ws.Range("L141").Value = 22500
$ws.Range("N141").Value = -32860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").Value = $null
$ws.Range("H132").Value = 1629.3103
$ws.Range("I132").Value = 1091.093
$ws.Range("K132").Value = 3273.279
$ws.Range("M132").Value = -743.2790000000005
